$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D and E (pushes batsman..sr from D..I to F..K)
$ws.Range("D1:E1").EntireColumn.Insert()

# Set new header values
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"

# Set new data row values
$ws.Range("D2").Value = "Mumbai Indians"
$ws.Range("E2").Value = "Chennai Super Kings"
